# Fruta / hortaliza, semanal
# Insert two new weekly rows of Papaya price data (Femacal de La Calera)
# at the top of the data block (row 94), shifting the existing rows
# 94-108 down to 96-110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 94 (pushes old rows 94-108 -> 96-110)
$ws.Range("A94:T95").EntireRow.Insert()

# --- New row 94: Papaya, Primera ---
$ws.Range("A94").Value = 3
$ws.Range("B94").Value = "Femacal de La Calera"
$ws.Range("C94").Value = "Coquimbo"
$ws.Range("D94").Value = 45180
$ws.Range("E94").Value = 5
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100108
$ws.Range("H94").Value = "Tropicales y subtropicales"
$ws.Range("I94").Value = 100108004
$ws.Range("J94").Value = "Papaya"
$ws.Range("K94").Value = "Cultivar IV Región"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 56
$ws.Range("N94").Value = 20000
$ws.Range("O94").Value = 20000
$ws.Range("P94").Value = 20000
$ws.Range("Q94").Value = '$/bandeja 10 kilos'
$ws.Range("R94").Value = "Provincia del Elquí"
$ws.Range("S94").Value = 2000
$ws.Range("T94").Value = 10

# --- New row 95: Papaya, Segunda ---
$ws.Range("A95").Value = 3
$ws.Range("B95").Value = "Femacal de La Calera"
$ws.Range("C95").Value = "Coquimbo"
$ws.Range("D95").Value = 45180
$ws.Range("E95").Value = 5
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100108
$ws.Range("H95").Value = "Tropicales y subtropicales"
$ws.Range("I95").Value = 100108004
$ws.Range("J95").Value = "Papaya"
$ws.Range("K95").Value = "Cultivar IV Región"
$ws.Range("L95").Value = "Segunda"
$ws.Range("M95").Value = 50
$ws.Range("N95").Value = 17000
$ws.Range("O95").Value = 17000
$ws.Range("P95").Value = 17000
$ws.Range("Q95").Value = '$/bandeja 10 kilos'
$ws.Range("R95").Value = "Provincia del Elquí"
$ws.Range("S95").Value = 1700
$ws.Range("T95").Value = 10
